$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings (e.g. "248.21")
# are stored as text, matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '41.603.08'
$ws.Range("E2").Value = '  -4.95%  '

# Row 3
$ws.Range("D3").Value = '2.234.52'
$ws.Range("E3").Value = '  -4.90%  '

# Row 4
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").Value = '248.21'
$ws.Range("E5").Value = '  +3.32%  '

# Row 6
$ws.Range("D6").Value = '0.631'
$ws.Range("E6").Value = '  -4.86%  '

# Row 7
$ws.Range("D7").Value = '71.40'
$ws.Range("E7").Value = '  -3.46%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  -5.42%  '

# Row 10
$ws.Range("D10").Value = '39.26'
$ws.Range("E10").Value = '  +6.44%  '

# Row 11
$ws.Range("D11").Value = '0.0963'
$ws.Range("E11").Value = '  -5.93%  '

# Row 12
$ws.Range("D12").Value = '58.35'
$ws.Range("E12").Value = '  -4.00%  '

# Row 13
$ws.Range("E13").Value = '  -3.68%  '

# Row 14
$ws.Range("D14").Value = '6.81'
$ws.Range("E14").Value = '  -6.07%  '

# Row 15
$ws.Range("D15").Value = '2.568.28'
$ws.Range("E15").Value = '  -4.95%  '

# Row 16
$ws.Range("D16").Value = '14.90'
$ws.Range("E16").Value = '  -8.54%  '

# Row 17
$ws.Range("D17").Value = '0.850'
$ws.Range("E17").Value = '  -8.67%  '

# Row 18
$ws.Range("D18").Value = '2.237.20'
$ws.Range("E18").Value = '  -5.81%  '

# Row 19
$ws.Range("D19").Value = '41.638.58'
$ws.Range("E19").Value = '  -4.72%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0964'
$ws.Range("E20").Value = '  -7.40%  '

# Row 21
$ws.Range("D21").Value = '73.10'
$ws.Range("E21").Value = '  -5.44%  '

# Row 22
$ws.Range("D22").Value = '6.14'
$ws.Range("E22").Value = '  -6.44%  '

# Row 23
$ws.Range("D23").Value = '233.83'
$ws.Range("E23").Value = '  -7.58%  '

# Row 24
$ws.Range("D24").Value = '2.08'
$ws.Range("E24").Value = '  +12.12%  '

# Row 25
$ws.Range("E25").Value = '  +0.03%  '

# Row 26
$ws.Range("D26").Value = '3.71'
$ws.Range("E26").Value = '  -0.88%  '

# Row 27
$ws.Range("D27").Value = '2.46'
$ws.Range("E27").Value = '  -1.44%  '

# Row 28
$ws.Range("D28").Value = '9.97'
$ws.Range("E28").Value = '  -6.50%  '

# Row 29
$ws.Range("D29").Value = '2.19'
$ws.Range("E29").Value = '  -4.36%  '

# Row 30
$ws.Range("D30").Value = '171.03'
$ws.Range("E30").Value = '  -2.65%  '

# Row 31
$ws.Range("D31").Value = '20.72'
$ws.Range("E31").Value = '  -6.32%  '

# Row 32
$ws.Range("D32").Value = '0.120'
$ws.Range("E32").Value = '  -6.65%  '

# Row 33
$ws.Range("D33").Value = '0.125'
$ws.Range("E33").Value = '  -6.97%  '

# Row 34
$ws.Range("D34").Value = '0.0717'
$ws.Range("E34").Value = '  -5.61%  '

# Row 35
$ws.Range("D35").Value = '5.30'
$ws.Range("E35").Value = '  -4.79%  '

# Row 36
$ws.Range("D36").Value = '4.67'
$ws.Range("E36").Value = '  -9.61%  '

# Row 37
$ws.Range("D37").Value = '4.05'
$ws.Range("E37").Value = '  +7.06%  '

# Row 38
$ws.Range("D38").Value = '25.02'
$ws.Range("E38").Value = '  +19.61%  '

# Row 39
$ws.Range("D39").Value = '0.0279'
$ws.Range("E39").Value = '  -0.78%  '

# Row 40
$ws.Range("D40").Value = '2.27'
$ws.Range("E40").Value = '  -4.41%  '

# Row 41
$ws.Range("D41").Value = '5.93'
$ws.Range("E41").Value = '  -10.46%  '

# Row 42
$ws.Range("D42").Value = '65.01'
$ws.Range("E42").Value = '  -1.87%  '

# Row 43
$ws.Range("D43").Value = '5.18'
$ws.Range("E43").Value = '  -4.71%  '

# Row 44
$ws.Range("D44").Value = '0.212'
$ws.Range("E44").Value = '  +5.71%  '

# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = '11.43'
$ws.Range("E45").Value = '  +14.23%  '

# Row 46
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '8.74'
$ws.Range("E46").Value = '  -4.42%  '

# Row 47
$ws.Range("D47").Value = '0.101'
$ws.Range("E47").Value = '  -7.10%  '

# Row 48
$ws.Range("B48").Value = 'SynthetixNetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D48").Value = '4.64'
$ws.Range("E48").Value = '  +4.98%  '

# Row 49
$ws.Range("B49").Value = 'BinanceUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.12%  '

# Row 50
$ws.Range("D50").Value = '1.18'
$ws.Range("E50").Value = '  -5.09%  '

# Row 51
$ws.Range("B51").Value = 'BitTorrent-New'
$ws.Range("C51").Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range("D51").Value = '0.0₃0150'
$ws.Range("E51").Value = '  +12.60%  '
